# The module number shown at the bottom of every slide comes from a
# textbox on the Slide Master (not an individual slide), so update it
# there.
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$shp = $master.Shapes.Item("TextBox 10")
$shp.TextFrame.TextRange.Text = "Module 1"
